# Edit the "Return of Biopsy Results within 1 Month" row of the RoR table:
#  1. Grow the row height (907 -> 1051 twips == 45.35 -> 52.55 points)
#  2. Append " past a month from biopsy" to the divisor text in column 2
#  3. Replace "Biopsy date" with "Past a month since biopsy" in column 3
#
# NOTE: Range.Find.Execute's Replace behaviour is NOT scoped to the Range
# it is invoked on in this runtime - it searches/replaces across the whole
# document regardless of which Range's Find object is used. Since the
# strings we need to change ("Biopsy date", "Number of biopsied
# participants") are repeated in many other rows of this table, we must
# NOT use Find/Replace here. Instead we compute exact character offsets
# for the target cells and assign Range.Text directly, which only touches
# the characters within that explicit Range.

$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

# --- Locate the target row by its first-cell label ---
$targetRow = $null
for ($i = 1; $i -le $table.Rows.Count; $i++) {
    $cellText = $table.Cell($i, 1).Range.Text
    if ($cellText -match "Return of Biopsy Results within 1 Month") {
        $targetRow = $table.Rows.Item($i)
        break
    }
}

# 1. Row height: 1051 twips / 20 = 52.55 points
$targetRow.Height = 52.55

# 2. Column 2: append " past a month from biopsy" after the existing
#    "Number of biopsied participants" suffix.
$cell2 = $table.Cell($targetRow.Index, 2)
$c2text = $cell2.Range.Text
$c2len = $c2text.Length - 2          # drop trailing CR + cell-mark chars
$c2start = $cell2.Range.Start
$oldSuffix = "Number of biopsied participants"
$newSuffix = "Number of biopsied participants past a month from biopsy"
if ($c2text.Substring(0, $c2len).EndsWith($oldSuffix)) {
    $suffixStart = $c2start + $c2len - $oldSuffix.Length
    $r2 = $d.Range($suffixStart, $c2start + $c2len)
    $r2.Text = $newSuffix
}

# 3. Column 3: replace the whole "Biopsy date" cell content with
#    "Past a month since biopsy".
$cell3 = $table.Cell($targetRow.Index, 3)
$c3text = $cell3.Range.Text
$c3len = $c3text.Length - 2          # drop trailing CR + cell-mark chars
$c3start = $cell3.Range.Start
$r3 = $d.Range($c3start, $c3start + $c3len)
$r3.Text = "Past a month since biopsy"
